$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing E2 value: "0.58510" -> "0.5851" (keep it stored as text) ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.5851"

# --- Add new data row 3 ---
# Make sure the whole new row starts out as "General" so the numeric rowid
# (column A) is written as a real number, matching row 2's A2 cell.
$ws.Range("A3:O3").NumberFormat = "General"

# These columns must hold text values that look numeric / date-like, so force
# a text format on them before assigning, otherwise Excel auto-converts them
# to number / date values.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3:J3").NumberFormat = "@"
$ws.Range("N3").NumberFormat = "@"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "BANKBILL"
$ws.Range("C3").Value = "AUD"
$ws.Range("D3").Value = "0.6433"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "2020-03-16"
$ws.Range("I3").Value = "null"
$ws.Range("J3").Value = "null"
$ws.Range("K3").Value = "003M"
$ws.Range("L3").Value = "DEPOSIT"
$ws.Range("M3").Value = "COMRLENDING"
$ws.Range("N3").Value = "null"
$ws.Range("O3").Value = "AUD,EUR"

$wb.Save()
